$d = $word.ActiveDocument

# NOTE: The diff's first hunk (", choose a process set size, ") only splits
# an existing run into several runs with byte-identical combined text
# (", choose a process se" + "t" + " size" + ", " == ", choose a process set size, ").
# There is no visible text/formatting change there, so it is intentionally
# left untouched.

# "... once you finish that trial (completed simulating all applications for
# the following process count), please press ..." -> "... process set) ..."
$d.Content.Find.Execute(
    " (completed simulating all applications for the following process count)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    " (completed simulating all applications for the following process set)", 2)

# "Trial 1: Process Count of 3" -> "Trial 1: Process Set Containing a Process Count of Three"
$d.Content.Find.Execute(
    "1: Process Count of 3",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "1: Process Set Containing a Process Count of Three", 2)

# "Trial 2: Process Count of 6" -> "Trial 2: Process Set Containing a Process Count of Six"
$d.Content.Find.Execute(
    "2: Process Count of 6",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "2: Process Set Containing a Process Count of Six", 2)

# "Trial 3: Process Count of 10 " (note trailing-space run that disappears)
# -> "Trial 3: Process Set Containing a Process Count of Ten"
$d.Content.Find.Execute(
    "3: Process Count of 10 ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "3: Process Set Containing a Process Count of Ten", 2)
